$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cryptocurrency price/volume text values scraped from coinranking.com.
# Some "Price" values look like plain decimal numbers to Excel (e.g. "213.25"),
# so we force those cells to Text format before assigning, then reset the style
# back to Normal so no stray formatting is left behind.

$ws.Range("D2").Value = "26.686.30"
$ws.Range("E2").Value = "  +1.55%  "
$ws.Range("D3").Value = "1.635.55"
$ws.Range("E3").Value = "  +1.02%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "213.25"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.506"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +4.10%  "
$ws.Range("E7").Value = "  -0.09%  "
$ws.Range("E8").Value = "  +2.71%  "
$ws.Range("E9").Value = "  +1.46%  "
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = "19.25"
$cell.Style = "Normal"
$ws.Range("E10").Value = "  +2.43%  "
$ws.Range("E11").Value = "  +3.21%  "
$ws.Range("D12").Value = "1.864.37"
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("D13").Value = "1.639.71"
$ws.Range("E13").Value = "  +1.29%  "
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("E15").Value = "  +1.94%  "
$ws.Range("D16").Value = "26.677.75"
$ws.Range("E16").Value = "  +1.49%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "63.33"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  +1.58%  "
$ws.Range("E18").Value = "  +2.47%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "219.01"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  +8.69%  "
$ws.Range("E20").Value = "  -0.08%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "4.32"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("E22").Value = "  +3.02%  "
$ws.Range("E23").Value = "  +1.26%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "1.93"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +0.99%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "148.34"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +3.09%  "
$ws.Range("E26").Value = "  -0.10%  "
$ws.Range("E27").Value = "  +1.63%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "6.92"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +5.79%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "15.53"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +2.47%  "
$ws.Range("E30").Value = "  -2.91%  "
$ws.Range("E31").Value = "  +0.17%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.33"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  +4.98%  "
$ws.Range("E33").Value = "  +1.92%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "1.52"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +1.36%  "
$ws.Range("E35").Value = "  -1.16%  "
$ws.Range("D36").Value = "1.197.63"
$ws.Range("E36").Value = "  +1.66%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.0174"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +5.75%  "
$ws.Range("E38").Value = "  +0.92%  "
$ws.Range("E39").Value = "  -0.14%  "
$ws.Range("E40").Value = "  +2.23%  "
$ws.Range("E41").Value = "  -1.53%  "
$ws.Range("E42").Value = "  +1.83%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.793"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Value = "1.772.85"
$ws.Range("E44").Value = "  +1.01%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "92.56"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -0.46%  "
$ws.Range("E46").Value = "  +2.65%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "54.90"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  +2.05%  "
$ws.Range("E48").Value = "  +0.91%  "
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "7.64"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  +5.36%  "
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.410"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  +0.57%  "
$ws.Range("E51").Value = "  +0.03%  "
